$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New employees appended below the existing table (rows 47-52), mirroring the
# layout/format of the row immediately above them (row 46).
$people = @(
    @{Row=47; A=49; B="LUIZ FERNANDO DE OLIVEITA CAETANO"; C="FERNANDINHO"},
    @{Row=48; A=50; B="JANILSON DOS SANTOS";               C="GEVÃO"},
    @{Row=49; A=51; B="GEOVAN DOS SANTOS";                 C="GEOVAN"},
    @{Row=50; A=52; B="ERCI CARLOS PEREIRA";                C="ERCI"},
    @{Row=51; A=53; B="EDUARDO DE ALMEIDA DOS SANTOS";      C="EDU"},
    @{Row=52; A=54; B="ANA CAROLINA DE ALBUQUERQUE PRADO";  C="CAROL"}
)

# Columns that use the workbook's "numeric" cell style (same style already
# used for these columns on every existing data row).
$styleCols = @("AH","AI","AQ","AW","AX","AY","BA","BC","BD","BE","BI","BJ","BL","BN")

# Values shared by every new row.
$sharedValues = @{
    AH = 0;    AI = 0;    AQ = 0.08
    AU = 220;  AV = 220
    AW = 0;    AX = 0;    AY = 0
    AZ = 0;    BA = 0;    BB = 0
    BC = 0;    BD = 0;    BE = 0
    BI = 0;    BJ = 0;    BK = 0
    BL = 0;    BN = 0;    BP = 0
}
$boolCols = @("BH", "BQ")

foreach ($p in $people) {
    $r = $p.Row
    $prev = $r - 1

    # Copy just the number-format style from the row above for the columns
    # that need it, so new cells match the sheet's existing formatting.
    foreach ($col in $styleCols) {
        $ws.Range($col + $prev).Copy() | Out-Null
        $ws.Range($col + $r).PasteSpecial(-4122) | Out-Null
    }

    $ws.Range("A$r").Value = $p.A
    $ws.Range("B$r").Value = $p.B
    $ws.Range("C$r").Value = $p.C

    foreach ($col in $sharedValues.Keys) {
        $ws.Range("$col$r").Value = $sharedValues[$col]
    }
    foreach ($col in $boolCols) {
        $ws.Range("$col$r").Value = $false
    }
}

# Grow the "Funcionários" named range to cover the newly added rows.
$nm = $wb.Names.Item(1)
$nm.RefersTo = "='Funcionários'!`$A`$1:`$BQ`$52"
